# Add a new row (30) of case data to the MuniEntryPleas sheet, mirroring
# the formatting of the row above it (row 29) and then filling in the
# new record's values (admin-code regex test row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MuniEntryPleas")

# Copy the formatting (styles) of the last existing data row onto the new
# row so the new cells pick up the same borders/fonts/number formats.
$ws.Range("A29:L29").Copy()
$ws.Range("A30:L30").PasteSpecial(-4122)

# Case / Sub Case / Lastname / Firstname / Charge / Code / Degree / Insurance
$ws.Range("A30").Value = "22TRC00570"
$ws.Range("B30").Value = "22TRC00570-A"
$ws.Range("C30").Value = "Kudela"
$ws.Range("D30").Value = "Justin"
$ws.Range("E30").Value = "TEST"
$ws.Range("F30").Value = "1501.17-5-04"
$ws.Range("G30").Value = "MM"
$ws.Range("H30").Value = "Y"

# Moving violation flag
$ws.Range("I30").Value = $true

# Atty Last / Atty First / Atty Type left blank for this record
$ws.Range("J30").Value = ""
$ws.Range("K30").Value = ""
$ws.Range("L30").Value = ""

# Leave the view with the newly-added row selected/visible.
$ws.Range("G30").Select()
